$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 758
$ws.Range("F6").Value = 2399
$ws.Range("F8").Value = 1774
$ws.Range("F9").Value = 3019
$ws.Range("F10").Value = 179
$ws.Range("F11").Value = 4486
$ws.Range("F12").Value = 396
$ws.Range("F13").Value = 219
$ws.Range("F15").Value = 565
$ws.Range("F16").Value = 266
$ws.Range("F17").Value = 620
$ws.Range("F18").Value = 233
$ws.Range("F20").Value = 111
$ws.Range("F21").Value = 312
$ws.Range("F22").Value = 4534
$ws.Range("F24").Value = 4032
$ws.Range("F25").Value = 1146
$ws.Range("F27").Value = 585
$ws.Range("F28").Value = 4380
$ws.Range("F29").Value = 93
$ws.Range("F30").Value = 624
$ws.Range("F31").Value = 603
$ws.Range("F32").Value = 571

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 28

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1042

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1042
$ws.Range("F8").Value = 758
$ws.Range("F9").Value = 2399
$ws.Range("F11").Value = 1774
$ws.Range("F13").Value = 3019
$ws.Range("F14").Value = 179
$ws.Range("F15").Value = 4486
$ws.Range("F16").Value = 396
$ws.Range("F17").Value = 219
$ws.Range("F19").Value = 565
$ws.Range("F20").Value = 266
$ws.Range("F21").Value = 620
$ws.Range("F22").Value = 233
$ws.Range("F23").Value = 28
$ws.Range("F25").Value = 111
$ws.Range("F26").Value = 312
$ws.Range("F27").Value = 4534
$ws.Range("F29").Value = 4034
$ws.Range("F30").Value = 1146
$ws.Range("F32").Value = 585
$ws.Range("F33").Value = 4380
$ws.Range("F34").Value = 93
$ws.Range("F35").Value = 624
$ws.Range("F36").Value = 603
$ws.Range("F37").Value = 571
